$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells in this sheet store text (inline strings in the original),
# even when the text looks numeric (e.g. "22.02") or has multiple dots
# (e.g. "26.861.18"). Force text format first so Excel COM does not
# auto-coerce numeric-looking strings into real numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.861.18'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -2.21%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.567.61'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.20%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -1.78%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '22.02'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.78%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.49%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0863'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.789.98'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.569.59'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.05%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.68%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '26.857.21'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -2.27%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.55'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -3.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '215.02'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.73%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.05'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.74'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '14.99'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.24%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.01%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.104'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0468'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.74%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.59%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.403.14'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.53%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.64%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.32%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.936'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.83%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0163'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.64%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.527'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.75%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.816'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.51%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.991'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.32'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.15%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.18'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.70%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '63.34'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.48%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.702.97'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '86.15'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0₇0986'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.41%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.35%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0492'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.97%  '
